# Update dynamic selection method to Signup Page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "LogIn_TestData"
$ws.Name = "LogIn_TestData"

# Increase the height of row 1 (header row)
$ws.Rows.Item(1).RowHeight = 33.6

# Update the active selection from R29 to B4
$ws.Range("B4").Select()
